$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 368 (shifts the existing row 368 and everything
# below it down by one row, e.g. old row 368 -> new row 369, ... old row 459 -> new row 460)
$ws.Rows(368).Insert()

# Populate the newly inserted row 368 with the new weekly price record
$ws.Range("A368").Value = 5
$ws.Range("B368").Value = "Macroferia Regional de Talca"
$ws.Range("C368").Value = "Maule"
$ws.Range("D368").Value = 44943
$ws.Range("E368").Value = 7
$ws.Range("F368").Value = 100114014
$ws.Range("G368").Value = "Betarraga"
$ws.Range("H368").Value = "Sin especificar"
$ws.Range("I368").Value = "Primera"
$ws.Range("J368").Value = 4000
$ws.Range("K368").Value = 650
$ws.Range("L368").Value = 700
$ws.Range("M368").Value = 675
$ws.Range("N368").Value = "$/paquete 5 unidades"
$ws.Range("O368").Value = "Región del Maule"
$ws.Range("P368").Value = 135
$ws.Range("Q368").Value = 5
$ws.Range("R368").Value = "Hortaliza"
